$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Trim the sheet down from A1:J6 to A1:C6 -------------------------------
# Remove every cell in columns D:J (for all used rows) so the sheet's
# dimension / used range shrinks to A1:C6, matching the target layout.
$ws.Range("D1:J6").Clear()

# --- Row 1 (no label in column A) ------------------------------------------
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2

# --- Row 2: carID -----------------------------------------------------------
$ws.Range("A2").Value = "carID"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2

# --- Row 3: speed2 -----------------------------------------------------------
$ws.Range("A3").Value = "speed2"
$ws.Range("B3").Value = 42.43
$ws.Range("C3").Value = 55.9

# --- Row 4: asma (values kept as text, not numbers) -------------------------
$ws.Range("A4").Value = "asma"
$ws.Range("B4").Value = "'21.22"
$ws.Range("C4").Value = "'59.72"

# --- Row 5: ceza_tutar --------------------------------------------------------
$ws.Range("A5").Value = "ceza_tutar"
$ws.Range("B5").Value = 1508.5
$ws.Range("C5").Value = 6440

# --- Row 6: hesaplanan_asma ---------------------------------------------------
$ws.Range("A6").Value = "hesaplanan_asma"
$ws.Range("B6").Value = 10
$ws.Range("C6").Value = 50
